# ajout ssr, had et psy
# Append a new data row (row 45) to the "Feuil1" worksheet describing an
# additional RHS field ("c" type, offset 181) labelled "ZAD".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New row of data
$ws.Range("B45").Value = "c"
$ws.Range("D45").Value = 181
$ws.Range("F45").Value = "ZAD"

# Update the view so the new row is visible / selected, mirroring the
# author's final cursor position after adding the row.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D46").Select()
